# Update "想去人数" (want-to-go count) figures in column F across sheets,
# matching the refreshed data snapshot referenced in the commit message.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 38
$ws1.Range("F3").Value = 104
$ws1.Range("F4").Value = 1509
$ws1.Range("F5").Value = 211
$ws1.Range("F7").Value = 287
$ws1.Range("F8").Value = 9893
$ws1.Range("F10").Value = 120
$ws1.Range("F14").Value = 6839
$ws1.Range("F15").Value = 1085
$ws1.Range("F16").Value = 633
$ws1.Range("F18").Value = 199

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 552

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 38
$ws4.Range("F3").Value = 104
$ws4.Range("F4").Value = 1509
$ws4.Range("F5").Value = 211
$ws4.Range("F8").Value = 287
$ws4.Range("F9").Value = 552
$ws4.Range("F11").Value = 9893
$ws4.Range("F13").Value = 120
$ws4.Range("F17").Value = 6839
$ws4.Range("F18").Value = 1085
$ws4.Range("F19").Value = 633
$ws4.Range("F21").Value = 199
